$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "N°1" list entry (row 12), re-using the label that used to be
# attached to E11 ("Partie 1"). Doing this before touching E11 makes the
# shared-string table grow in the same order as in the authoritative edit.
$ws.Range("B11").Copy()
$ws.Range("B12").PasteSpecial(-4122)
$ws.Range("B12").Value = 44185

$ws.Range("C11").Copy()
$ws.Range("C12").PasteSpecial(-4122)
$ws.Range("C12").Value = "Anthony"

$ws.Range("D12").Value = 6

$ws.Range("E12").Value = "N°1"

# The old E11 entry now reads "Tout" instead of "Partie 1".
$ws.Range("E11").Value = "Tout"

$ws.Range("E12").Select()
